$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# Row 2: Amapá
$ws.Range("A2").Value = "Amapá"
$ws.Range("B2").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C2").Value = 0.02051623430250304
$ws.Range("D2").Value = "1º"

# Row 3: Mato Grosso do Sul
$ws.Range("A3").Value = "Mato Grosso do Sul"
$ws.Range("B3").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C3").Value = 0.005608039820102051
$ws.Range("D3").Value = "2º"

# Row 4: Roraima
$ws.Range("A4").Value = "Roraima"
$ws.Range("B4").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C4").Value = 0.003111664295874839
$ws.Range("D4").Value = "3º"

# Row 5: Ceará
$ws.Range("A5").Value = "Ceará"
$ws.Range("B5").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C5").Value = 0.00047684956933014
$ws.Range("D5").Value = "4º"

# Row 6: Piauí
$ws.Range("A6").Value = "Piauí"
$ws.Range("B6").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C6").Value = -0.00001699323063175662
$ws.Range("D6").Value = "5º"

# Row 7: Bahia
$ws.Range("A7").Value = "Bahia"
$ws.Range("B7").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C7").Value = -0.001434825472076906
$ws.Range("D7").Value = "6º"

# Row 8: Sergipe (A unchanged)
$ws.Range("B8").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C8").Value = -0.002622913429601859
$ws.Range("D8").Value = "10º"

# Row 9: Nordeste (A unchanged)
$ws.Range("B9").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C9").Value = -0.005027825232064465

# Row 10: Brasil (A unchanged)
$ws.Range("B10").Value = "Diferença 2023/03 - 2022/03"
$ws.Range("C10").Value = -0.006086099171784931
